$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting note row for 29/7 ("buoi hop ngay 29/7")
$ws.Range("A19").Value = 45502
$ws.Range("A19").NumberFormat = "mm-dd-yy"

$ws.Range("B19").Value = "Erd + class : diagram "

# Put the cursor where the user left it after typing the new row
$ws.Range("C19").Select() | Out-Null
